$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 21
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = $null
$ws.Range("H54").Value = 9333.333000000001
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = $null
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = $null
$ws.Range("M54").Value = -9514
$ws.Range("N54").Value = -8972
$ws.Range("H74").Value = 7999.5
$ws.Range("I74").Value = 7999.5
$ws.Range("K74").Value = 7999.5
$ws.Range("M74").Value = -7063.5
$ws.Range("H77").Value = 7999.5
$ws.Range("I77").Value = 7999.5
$ws.Range("K77").Value = 39997.5
$ws.Range("M77").Value = -35317.5
$ws.Range("H96").Value = 803.1429000000001
$ws.Range("I96").Value = 710.4
$ws.Range("J96").Value = 854.6667
$ws.Range("K96").Value = 2131.2
$ws.Range("L96").Value = 2564.0001
$ws.Range("M96").Value = -758.1999999999998
$ws.Range("N96").Value = -5310.0001
$ws.Range("H108").Value = 47500
$ws.Range("J108").Value = 47500
$ws.Range("L108").Value = 47500
$ws.Range("N108").Value = -55180

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1083.875
$ws.Range("I74").Value = 778.5
$ws.Range("K74").Value = 778.5
$ws.Range("M74").Value = 95.5
$ws.Range("H77").Value = 1083.875
$ws.Range("I77").Value = 778.5
$ws.Range("K77").Value = 3892.5
$ws.Range("M77").Value = 475.5
$ws.Range("H97").Value = 55557960
$ws.Range("H102").Value = 19183004
$ws.Range("I102").Value = 1101304.9
$ws.Range("K102").Value = 1101304.9
$ws.Range("M102").Value = -1099682.9
$ws.Range("H106").Value = 9995.5
$ws.Range("J106").Value = 9995.5
$ws.Range("L106").Value = 9995.5
$ws.Range("N106").Value = -12519.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1993
$ws.Range("I86").Value = 2624
$ws.Range("J86").Value = 100
$ws.Range("K86").Value = 2624
$ws.Range("L86").Value = 100
$ws.Range("M86").Value = -1501
$ws.Range("N86").Value = -2346
$ws.Range("H89").Value = 1993
$ws.Range("I89").Value = 2624
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 13120
$ws.Range("L89").Value = 500
$ws.Range("M89").Value = -7504
$ws.Range("N89").Value = -11732
$ws.Range("H94").Value = 221601.8
$ws.Range("I94").Value = 221601.8
$ws.Range("K94").Value = 221601.8
$ws.Range("M94").Value = -221150.8
$ws.Range("H117").Value = 49500
$ws.Range("J117").Value = 49500
$ws.Range("L117").Value = 49500
$ws.Range("N117").Value = -58678

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4309
$ws.Range("J10").Value = 7499
$ws.Range("L10").Value = 7499
$ws.Range("N10").Value = -7777
$ws.Range("H22").Value = 659.46155
$ws.Range("I22").Value = 726.4545000000001
$ws.Range("J22").Value = 291
$ws.Range("K22").Value = 726.4545000000001
$ws.Range("L22").Value = 291
$ws.Range("M22").Value = -376.4545000000001
$ws.Range("N22").Value = -991
$ws.Range("H107").Value = 917.3
$ws.Range("I107").Value = 824.8570999999999
$ws.Range("K107").Value = 824.8570999999999
$ws.Range("M107").Value = 1095.1429
$ws.Range("H134").Value = 1287.4286
$ws.Range("I134").Value = 874.5
$ws.Range("J134").Value = 1838
$ws.Range("K134").Value = 2623.5
$ws.Range("L134").Value = 5514
$ws.Range("M134").Value = -88.5
$ws.Range("N134").Value = -10584

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 345.14285
$ws.Range("I2").Value = 367.84616
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 2207.07696
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -2094.07696
$ws.Range("N2").Value = -526
$ws.Range("H129").Value = 555.8
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null
$ws.Range("H131").Value = 1412.7778
$ws.Range("I131").Value = 964.375
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 2893.125
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 2146.875
$ws.Range("N131").Value = -25080
$ws.Range("H132").Value = 950.4286
$ws.Range("J132").Value = 891
$ws.Range("L132").Value = 8019
$ws.Range("N132").Value = -13079
$ws.Range("H137").Value = 999
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 2103
$ws.Range("N137").Value = $null
$ws.Range("H141").Value = 9972.5
$ws.Range("I141").Value = 9965
$ws.Range("K141").Value = 29895
$ws.Range("M141").Value = -24715

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3110
$ws.Range("I80").Value = 3320
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 3320
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -2322
$ws.Range("N80").Value = -4896
$ws.Range("H83").Value = 3110
$ws.Range("I83").Value = 3320
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 16600
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -11608
$ws.Range("N83").Value = -24484

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2951.818
$ws.Range("I22").Value = 2909.0908
$ws.Range("J22").Value = 2994.5454
$ws.Range("K22").Value = 2909.0908
$ws.Range("L22").Value = 2994.5454
$ws.Range("M22").Value = -2614.0908
$ws.Range("N22").Value = -3584.5454
$ws.Range("H27").Value = 2951.818
$ws.Range("I27").Value = 2909.0908
$ws.Range("J27").Value = 2994.5454
$ws.Range("K27").Value = 2909.0908
$ws.Range("L27").Value = 2994.5454
$ws.Range("M27").Value = -2802.0908
$ws.Range("N27").Value = -3208.5454
$ws.Range("H46").Value = 1395
$ws.Range("I46").Value = 1375
$ws.Range("J46").Value = 1408.3334
$ws.Range("K46").Value = 1375
$ws.Range("L46").Value = 1408.3334
$ws.Range("M46").Value = -1187
$ws.Range("N46").Value = -1784.3334
$ws.Range("H93").Value = 55556404
$ws.Range("I93").Value = 66667490
$ws.Range("K93").Value = 66667490
$ws.Range("M93").Value = -66666242

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 17000
$ws.Range("J28").Value = 17000
$ws.Range("L28").Value = 17000
$ws.Range("N28").Value = -17696
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H122").Value = 1582.9231
$ws.Range("I122").Value = 1186.5555
$ws.Range("K122").Value = 3559.6665
$ws.Range("M122").Value = -1109.6665
